# The date line "3/3/2019" must become "10/3/2019".
# In the true editing session, the author selected just the first "3" and
# typed "10" over it; Word's "_GoBack" last-edit bookmark then tracks that
# edit location, landing right after the newly typed "10" (i.e. between
# the first run and the rest of the date), while the remaining runs
# ("/", "3", "/2019") are left untouched.

$d = $word.ActiveDocument

# Locate the "Namen en datum" date paragraph by finding the literal date text.
$found = $d.Content.Find
$found.ClearFormatting()
$found.Text = "3/3/2019"
$found.Execute() | Out-Null

$dateRange = $found.Parent
$start = $dateRange.Start

# Range covering just the first character ("3") of the date.
$firstCharRange = $d.Range($start, $start + 1)

# Re-seat the "_GoBack" bookmark to sit right after that first character,
# matching where Word leaves it after an in-place edit there.
$bmRange = $d.Range($start + 1, $start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Now replace the first character's text ("3" -> "10"); this run keeps its
# own formatting/run and the later runs stay split apart by the bookmark.
$firstCharRange.Text = "10"
